# TimeLog.xlsx - "ref timelog for notes"
#
# Adds a new day entry (row 14) to the Sheet1 log table:
#   B14 - date (14-Mar-2024, serial 45365), formatted like the other date cells
#   C14 - time spent (3 hours)
#   D14 - "What Was Achieved" note
#   E14 - "Notes" note
# and moves the active-cell selection down, the way a user would after
# finishing typing the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- new log row -----------------------------------------------------
$ws.Range("B14").Value = 45365
$ws.Range("B14").NumberFormat = "d-mmm"
$ws.Range("B14").VerticalAlignment = -4107   # xlBottom

$ws.Range("C14").Value = 3

$ws.Range("D14").Value = "Added a bounding box for map. Fixed prefab sizes and hit boxes. Created dictionary for all tiles"
$ws.Range("E14").Value = "Dictionary created. Need to create one for finding out where player ship placements are and then can link to attacking script."

# Row 14 holds longer wrapped text, so it ends up taller than the others.
$ws.Rows.Item(14).RowHeight = 28

# --- recalc formulas (Total Time Taken / Time Left) -------------------
$wb.Application.Calculate()

# --- leave the selection where the user would land afterwards ---------
$ws.Range("D17").Select()
